$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1166.1111
$ws.Range("J29").Value = 150
$ws.Range("L29").Value = 450
$ws.Range("N29").Value = -1012

$ws.Range("H38").Value = 1908.8
$ws.Range("I38").Value = 1386.25
$ws.Range("K38").Value = 4158.75
$ws.Range("M38").Value = -3786.75

$ws.Range("H43").Value = 2413.4285
$ws.Range("I43").Value = 1347.5
$ws.Range("J43").Value = 2839.8
$ws.Range("K43").Value = 1347.5
$ws.Range("L43").Value = 2839.8
$ws.Range("M43").Value = -1278.5
$ws.Range("N43").Value = -2977.8

$ws.Range("H49").Value = 255
$ws.Range("J49").Value = 200
$ws.Range("L49").Value = 600
$ws.Range("N49").Value = -872

$ws.Range("H62").Value = 8501.866
$ws.Range("I62").Value = 6857.364
$ws.Range("K62").Value = 6857.364
$ws.Range("M62").Value = -6233.364

$ws.Range("H65").Value = 8501.866
$ws.Range("I65").Value = 6857.364
$ws.Range("K65").Value = 34286.82
$ws.Range("M65").Value = -31166.82

$ws.Range("H100").Value = 2055.9375
$ws.Range("I100").Value = 1990.8334
$ws.Range("K100").Value = 1990.8334
$ws.Range("M100").Value = -1449.8334

$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139

$ws.Range("H137").Value = 34490492
$ws.Range("I137").Value = 100001290
$ws.Range("J137").Value = 11125.947
$ws.Range("K137").Value = 300003870
$ws.Range("L137").Value = 33377.841
$ws.Range("M137").Value = -300001320
$ws.Range("N137").Value = -38477.841

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5099.6
$ws.Range("I63").Value = 4999.3335
$ws.Range("J63").Value = 5250
$ws.Range("K63").Value = 4999.3335
$ws.Range("L63").Value = 5250
$ws.Range("M63").Value = -4313.3335
$ws.Range("N63").Value = -6622

$ws.Range("H66").Value = 5099.6
$ws.Range("I66").Value = 4999.3335
$ws.Range("J66").Value = 5250
$ws.Range("K66").Value = 24996.6675
$ws.Range("L66").Value = 26250
$ws.Range("M66").Value = -21564.6675
$ws.Range("N66").Value = -33114

$ws.Range("H74").Value = 3482169
$ws.Range("I74").Value = 7939458
$ws.Range("J74").Value = 15388.889
$ws.Range("K74").Value = 7939458
$ws.Range("L74").Value = 15388.889
$ws.Range("M74").Value = -7938584
$ws.Range("N74").Value = -17136.889

$ws.Range("H77").Value = 3482169
$ws.Range("I77").Value = 7939458
$ws.Range("J77").Value = 15388.889
$ws.Range("K77").Value = 39697290
$ws.Range("L77").Value = 76944.44499999999
$ws.Range("M77").Value = -39692922
$ws.Range("N77").Value = -85680.44499999999

$ws.Range("H132").Value = 760529.3
$ws.Range("I132").Value = 895863.5600000001
$ws.Range("K132").Value = 2687590.68
$ws.Range("M132").Value = -2685060.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 32471.5
$ws.Range("J30").Value = 32471.5
$ws.Range("L30").Value = 32471.5
$ws.Range("N30").Value = -32721.5

$ws.Range("H86").Value = 2004.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2004.5
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 2004.5
$ws.Range("N86").Value = -4250.5

$ws.Range("H89").Value = 2004.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2004.5
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 10022.5
$ws.Range("N89").Value = -21254.5

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H99").Value = 7517.1177
$ws.Range("I99").Value = 10581
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 10581
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -9083
$ws.Range("N99").Value = -4896

$ws.Range("H107").Value = 1550
$ws.Range("I107").Value = 1057.1428
$ws.Range("K107").Value = 1057.1428
$ws.Range("M107").Value = 862.8571999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 921.41174
$ws.Range("I19").Value = 354.0625
$ws.Range("K19").Value = 354.0625
$ws.Range("M19").Value = -184.0625

$ws.Range("H24").Value = 921.41174
$ws.Range("I24").Value = 354.0625
$ws.Range("K24").Value = 354.0625
$ws.Range("M24").Value = -184.0625

$ws.Range("H32").Value = 1975
$ws.Range("I32").Value = 1975
$ws.Range("K32").Value = 1975
$ws.Range("M32").Value = -1659

$ws.Range("H107").Value = 354.5
$ws.Range("J107").Value = 373
$ws.Range("L107").Value = 373
$ws.Range("N107").Value = -4213

$ws.Range("H132").Value = 4302.5713
$ws.Range("I132").Value = 4249.077
$ws.Range("K132").Value = 12747.231
$ws.Range("M132").Value = -10217.231

$ws.Range("H141").Value = 158856.4
$ws.Range("J141").Value = 183612.14
$ws.Range("L141").Value = 183612.14
$ws.Range("N141").Value = -193972.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19713858
$ws.Range("I4").Value = 30606750
$ws.Range("J4").Value = 2909.9524
$ws.Range("K4").Value = 91820250
$ws.Range("L4").Value = 8729.8572
$ws.Range("M4").Value = -91820138
$ws.Range("N4").Value = -8953.8572

$ws.Range("H34").Value = 443.33334
$ws.Range("J34").Value = 650
$ws.Range("L34").Value = 1950
$ws.Range("N34").Value = -2118

$ws.Range("H39").Value = 3840
$ws.Range("J39").Value = 6000
$ws.Range("L39").Value = 18000
$ws.Range("N39").Value = -18588

$ws.Range("H132").Value = 1381.5264
$ws.Range("I132").Value = 979.2
$ws.Range("J132").Value = 1828.5555
$ws.Range("K132").Value = 8812.800000000001
$ws.Range("L132").Value = 16456.9995
$ws.Range("M132").Value = -6282.800000000001
$ws.Range("N132").Value = -21516.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 74804.25
$ws.Range("I122").Value = 104533.63
$ws.Range("J122").Value = 9399.6
$ws.Range("K122").Value = 313600.89
$ws.Range("L122").Value = 28198.8
$ws.Range("M122").Value = -311150.89
$ws.Range("N122").Value = -33098.8

$ws.Range("H126").Value = 7618.077
$ws.Range("I126").Value = 10458.25
$ws.Range("J126").Value = 3073.8
$ws.Range("K126").Value = 31374.75
$ws.Range("L126").Value = 9221.400000000001
$ws.Range("M126").Value = -28904.75
$ws.Range("N126").Value = -14161.4

$ws.Range("H132").Value = 31155.54
$ws.Range("I132").Value = 30001.625
$ws.Range("J132").Value = 33001.8
$ws.Range("K132").Value = 90004.875
$ws.Range("L132").Value = 99005.40000000001
$ws.Range("M132").Value = -87474.875
$ws.Range("N132").Value = -104065.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 29000
$ws.Range("J38").Value = 29000
$ws.Range("L38").Value = 29000
$ws.Range("N38").Value = -29820

$ws.Range("H40").Value = 6945
$ws.Range("I40").Value = 6416.8335
$ws.Range("K40").Value = 6416.8335
$ws.Range("M40").Value = -6280.8335

$ws.Range("H61").Value = 11288.223
$ws.Range("I61").Value = 10792.667
$ws.Range("K61").Value = 10792.667
$ws.Range("M61").Value = -10590.667

$ws.Range("H68").Value = 1827.6666
$ws.Range("I68").Value = 1817.5883
$ws.Range("K68").Value = 1817.5883
$ws.Range("M68").Value = -1068.5883

$ws.Range("H71").Value = 1827.6666
$ws.Range("I71").Value = 1817.5883
$ws.Range("K71").Value = 9087.941499999999
$ws.Range("M71").Value = -5343.941499999999

$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -459

$ws.Range("H113").Value = 11288.223
$ws.Range("I113").Value = 10792.667
$ws.Range("K113").Value = 10792.667
$ws.Range("M113").Value = -8622.666999999999

$ws.Range("H136").Value = 6762498.5
$ws.Range("I136").Value = 9617894
$ws.Range("K136").Value = 28853682
$ws.Range("M136").Value = -28851132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 26035.25
$ws.Range("J52").Value = 26035.25
$ws.Range("L52").Value = 26035.25
$ws.Range("N52").Value = -26487.25

$ws.Range("H107").Value = 1282.55
$ws.Range("I107").Value = 425.5
$ws.Range("J107").Value = 4710.75
$ws.Range("K107").Value = 1276.5
$ws.Range("L107").Value = 14132.25
$ws.Range("M107").Value = 643.5
$ws.Range("N107").Value = -17972.25

$ws.Range("H135").Value = 102499.5
$ws.Range("J135").Value = 102499.5
$ws.Range("L135").Value = 102499.5
$ws.Range("N135").Value = -112639.5

$ws.Range("H136").Value = 12816171
$ws.Range("I136").Value = 2718803
$ws.Range("J136").Value = 66668800
$ws.Range("K136").Value = 8156409
$ws.Range("L136").Value = 200006400
$ws.Range("M136").Value = -8153859
$ws.Range("N136").Value = -200011500
